$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 559, shifting existing rows 559-677 down to 560-678.
$ws.Rows.Item(559).Insert()

# Populate the newly inserted row 559 with the new data record.
$ws.Cells.Item(559, 1).Value = 9
$ws.Cells.Item(559, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(559, 3).Value = "Metropolitana"
$ws.Cells.Item(559, 4).Value = 45209
$ws.Cells.Item(559, 5).Value = 13
$ws.Cells.Item(559, 6).Value = 100112012
$ws.Cells.Item(559, 7).Value = "Espinaca"
$ws.Cells.Item(559, 8).Value = "Sin especificar"
$ws.Cells.Item(559, 9).Value = "Primera"
$ws.Cells.Item(559, 10).Value = 160
$ws.Cells.Item(559, 11).Value = 7000
$ws.Cells.Item(559, 12).Value = 8000
$ws.Cells.Item(559, 13).Value = 7500
$ws.Cells.Item(559, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(559, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(559, 16).Value = 750
$ws.Cells.Item(559, 17).Value = 10
$ws.Cells.Item(559, 18).Value = "Hortaliza"
